# Apply the "Changed Designator K4 to Q1 and updated 'Stückliste'" edit.
#
# This adds Designator values (Position column, A) to rows 4-11 and fills
# in the previously-empty rows 10 and 11 with a new LED part and a new
# MOSFET part, respectively. Dependent formulas (I10, I11, I32) recalc
# automatically.
#
# Cell values are written in the same order the original author appears to
# have entered them (new parts' detail columns first, designators last) so
# that newly interned shared strings line up with the authoritative file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: new LED part ---------------------------------------------------
$ws.Range("C10").Value = "APTD3216SRCPRV; SMD RED"
$ws.Range("B10").Value = "LED "

# --- Row 11: new MOSFET part ------------------------------------------------
$ws.Range("F11").Value = "863-NTR4501NT1G "
$ws.Range("D11").Value = "SOT-23-3"
$ws.Range("C11").Value = "NTR4501NT1G; 20V; 3.2A; N-Kanal"
$ws.Range("B11").Value = "MOSFET"

# Replace the row-10 order code (formerly had stray leading space)
$ws.Range("F10").Value = "604-APTD3216SRCPRV "

# --- Designators (column A, Position) --------------------------------------
$ws.Range("A4").Value  = "K3"
$ws.Range("A5").Value  = "K3"
$ws.Range("A6").Value  = "K5"
$ws.Range("A7").Value  = "K1"
$ws.Range("A8").Value  = "X5"
$ws.Range("A9").Value  = "K1"
$ws.Range("A10").Value = "P1"
$ws.Range("A11").Value = "Q1"

# --- Remaining (non-string) cells of the two new rows -----------------------
$ws.Range("D10").Value = 1206
$ws.Range("E10").Value = "Mouser"
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0.356

$ws.Range("E11").Value = "Mouser"
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 0.338

# --- Row heights (both rows now wrap onto two lines) ------------------------
$ws.Rows.Item(10).RowHeight = 33
$ws.Rows.Item(11).RowHeight = 33

# --- View state: scroll down a bit, zoom to 100% and select A12 ------------
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("A12").Select()
